$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- Update the two existing accrual rows just above the new 2024 block ---
$ws.Range("C51").Value = 1.25
$ws.Range("C52").Value = 1.25

# --- Insert a new row before row 53 (everything below shifts down by one) ---
$ws.Rows("53:53").Insert()

# Make the new row look like the other year-separator rows (10, 14, 27, 40)
$ws.Range("A10:K10").Copy()
$ws.Range("A53:K53").PasteSpecial(-4122)

# New "2024" year-separator row
$ws.Range("A53").Value = "'2024"
$ws.Range("G53").Formula = '=IF(ISBLANK([@EARNED]),"",[@EARNED])'

# The row that used to be 53 (01/01/2024 accrual) is now row 54 - fill in its data
$ws.Range("B54").Value = "VL(6-0-0)"
$ws.Range("D54").Value = 6
$ws.Range("K54").Value = "01/05,08-12/2024"

# --- Expand the table (Table15) to cover the newly inserted row ---
$lo = $ws.ListObjects.Item("Table15")
$lo.Resize($ws.Range("A8:K100"))
